$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 8.351570000000001
$ws.Cells.Item(2, 8).Value = 25.05471
$ws.Cells.Item(2, 9).Value = 0.3629556103554933
$ws.Cells.Item(2, 10).Value = 0.3629556103554933
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 1.707786666666667
$ws.Cells.Item(2, 14).Value = 5.12336
$ws.Cells.Item(2, 15).Value = 0.8764025646701329
$ws.Cells.Item(2, 16).Value = 0.8764025646701328
$ws.Cells.Item(2, 17).Value = 14.26269989173333
$ws.Cells.Item(2, 18).Value = 128.3642990256
$ws.Cells.Item(2, 19).Value = 0.3180952277769678
$ws.Cells.Item(2, 20).Value = 0.3180952277769677

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 8.351570000000001
$ws.Cells.Item(3, 8).Value = 25.05471
$ws.Cells.Item(3, 9).Value = 0.3629556103554933
$ws.Cells.Item(3, 10).Value = 0.3629556103554933
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 0.240846
$ws.Cells.Item(3, 14).Value = 0.7225379999999999
$ws.Cells.Item(3, 15).Value = 0.1235974353298672
$ws.Cells.Item(3, 16).Value = 0.1235974353298672
$ws.Cells.Item(3, 17).Value = 2.01144222822
$ws.Cells.Item(3, 18).Value = 18.10298005398
$ws.Cells.Item(3, 19).Value = 0.04486038257852556
$ws.Cells.Item(3, 20).Value = 0.04486038257852556

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 11.216696
$ws.Cells.Item(4, 8).Value = 33.650088
$ws.Cells.Item(4, 9).Value = 0.4874727437897329
$ws.Cells.Item(4, 10).Value = 0.487472743789733
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 1.707786666666667
$ws.Cells.Item(4, 14).Value = 5.12336
$ws.Cells.Item(4, 15).Value = 0.8764025646701329
$ws.Cells.Item(4, 16).Value = 0.8764025646701328
$ws.Cells.Item(4, 17).Value = 19.15572387285333
$ws.Cells.Item(4, 18).Value = 172.40151485568
$ws.Cells.Item(4, 19).Value = 0.4272223628641085
$ws.Cells.Item(4, 20).Value = 0.4272223628641085

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 11.216696
$ws.Cells.Item(5, 8).Value = 33.650088
$ws.Cells.Item(5, 9).Value = 0.4874727437897329
$ws.Cells.Item(5, 10).Value = 0.487472743789733
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 0.240846
$ws.Cells.Item(5, 14).Value = 0.7225379999999999
$ws.Cells.Item(5, 15).Value = 0.1235974353298672
$ws.Cells.Item(5, 16).Value = 0.1235974353298672
$ws.Cells.Item(5, 17).Value = 2.701496364815999
$ws.Cells.Item(5, 18).Value = 24.31346728334399
$ws.Cells.Item(5, 19).Value = 0.06025038092562444
$ws.Cells.Item(5, 20).Value = 0.06025038092562444

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 3.441627666666667
$ws.Cells.Item(6, 8).Value = 10.324883
$ws.Cells.Item(6, 9).Value = 0.1495716458547737
$ws.Cells.Item(6, 10).Value = 0.1495716458547737
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 1.707786666666667
$ws.Cells.Item(6, 14).Value = 5.12336
$ws.Cells.Item(6, 15).Value = 0.8764025646701329
$ws.Cells.Item(6, 16).Value = 0.8764025646701328
$ws.Cells.Item(6, 17).Value = 5.877565840764444
$ws.Cells.Item(6, 18).Value = 52.89809256688
$ws.Cells.Item(6, 19).Value = 0.1310849740290565
$ws.Cells.Item(6, 20).Value = 0.1310849740290565

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 3.441627666666667
$ws.Cells.Item(7, 8).Value = 10.324883
$ws.Cells.Item(7, 9).Value = 0.1495716458547737
$ws.Cells.Item(7, 10).Value = 0.1495716458547737
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 0.240846
$ws.Cells.Item(7, 14).Value = 0.7225379999999999
$ws.Cells.Item(7, 15).Value = 0.1235974353298672
$ws.Cells.Item(7, 16).Value = 0.1235974353298672
$ws.Cells.Item(7, 17).Value = 0.828902257006
$ws.Cells.Item(7, 18).Value = 7.460120313053999
$ws.Cells.Item(7, 19).Value = 0.01848667182571719
$ws.Cells.Item(7, 20).Value = 0.01848667182571719
